$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Change 1: "Creación de clases, métodos y atributos" -> split "clases"
# into its own run and apply strike-through to it.
# -----------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("Creación de clases, métodos y atributos", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $paraStart1 = $r1.Start
    $paraEnd1 = $r1.End
    $scoped1 = $d.Range($paraStart1, $paraEnd1)
    $scoped1.Find.Execute("clases", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $target1 = $d.Range($scoped1.Start, $scoped1.End)
    $target1.Font.StrikeThrough = 1
}

# -----------------------------------------------------------------------
# Change 2 & 3: "Rodar sin deslizar" and "Fricción al frenar" paragraphs
# get full strike-through (paragraph mark + run).
# -----------------------------------------------------------------------
foreach ($txt in @("Rodar sin deslizar", "Fricción al frenar")) {
    $rr = $d.Content
    $f = $rr.Find.Execute($txt, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($f) {
        $pg = $d.Range($rr.Start, $rr.End)
        $pgParagraph = $pg.Paragraphs.Item(1)
        $pgParagraph.Range.Font.StrikeThrough = 1
    }
}

# -----------------------------------------------------------------------
# Change 4: "Creación de 3 mundos." -> "Creación de 3 mundos:" (two runs)
# plus three new sub-bullets ("1er mundo", "2ndo mundo", "3er mundo")
# and a trailing blank ListParagraph (no numbering) at ind left=1440.
# -----------------------------------------------------------------------
$r4 = $d.Content
$found4 = $r4.Find.Execute("mundos.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found4) {
    $target4 = $d.Range($r4.Start, $r4.End)
    $xml4 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-CO"/></w:rPr><w:t>mundos</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-CO"/></w:rPr><w:t>:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target4.InsertXML($xml4)
}

$r4b = $d.Content
$found4b = $r4b.Find.Execute("Creación de 3 mundos:")
if ($found4b) {
    $endOfPara4 = $r4b.End
    $insPoint4 = $d.Range($endOfPara4, $endOfPara4)
    $xml4b = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
        '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="5"/></w:numPr><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-CO"/></w:rPr><w:t>1er mundo</w:t></w:r></w:p>' +
        '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="5"/></w:numPr><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-CO"/></w:rPr><w:t>2ndo mundo</w:t></w:r></w:p>' +
        '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="5"/></w:numPr><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-CO"/></w:rPr><w:t>3er mundo</w:t></w:r></w:p>' +
        '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1440"/><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-CO"/></w:rPr></w:pPr></w:p>' +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $insPoint4.InsertXML($xml4b)
}

# -----------------------------------------------------------------------
# Change 5: add <w:lastRenderedPageBreak/> before "Implementación 3
# niveles de dificultad".
# -----------------------------------------------------------------------
$r5 = $d.Content
$found5 = $r5.Find.Execute("Implementación 3 niveles de dificultad")
if ($found5) {
    $target5 = $d.Range($r5.Start, $r5.End)
    $xml5 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-CO"/></w:rPr><w:lastRenderedPageBreak/><w:t>Implementación 3 niveles de dificultad</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target5.InsertXML($xml5)
}
